$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "forum of trajan"
$ws.Range("A18").Value = "hadrian bust"
$ws.Range("A19").Value = "panteon"
$ws.Range("A20").Value = "column of marcus"

$ws.Range("B17").Value = "rome"
$ws.Range("B18").Value = "rome"
$ws.Range("B19").Value = "rome"
$ws.Range("B20").Value = "rome"

$ws.Range("C17").Value = "112 ce"
$ws.Range("C18").Value = "117-120 ce"
$ws.Range("C19").Value = "118-125 ce"
$ws.Range("C20").Value = "180-192 ce"

$ws.Range("D17").Value = "marble"
$ws.Range("D18").Value = "marble"
$ws.Range("D19").Value = "marble"
$ws.Range("D20").Value = "marble"

$ws.Range("E17").Value = "quiz2/17.png"
$ws.Range("E18").Value = "quiz2/18.png"
$ws.Range("E19").Value = "quiz2/19.png"
$ws.Range("E20").Value = "quiz2/20.png"

$ws.Range("A12").Select()
$excel.ActiveWindow.Zoom = 145
